# "Same context, fixed issue 12"
#
# - Issue #12 (row 14 on the "Issues" sheet) gets its Root Cause / Solution
#   filled in and its Status flipped from "Open" to "Fixed".
# - The active/selected sheet moves from "Features" back to "Issues", with
#   the selection landing on F15.

$wb = $excel.ActiveWorkbook

$issues = $wb.Worksheets.Item("Issues")

# Fill in Root Cause (D14) and Solution (E14) for issue #12, and mark it Fixed.
$issues.Range("D14").Value = "在从DB取出m的时候，用了多个context。例如某个battery,在Batteries = dbContext.Batteries…的时候加载了一次，在programs.subprograms.testrecords.assignedbattery的时候又加载了一次。而这两次加载位于不同的context，所以即使他们的id相同，但他们其实是不同的instance。于是，BatteryVM绑定了model1的事件，而ProgramVM使得model2发出事件,并不会通知BatteryVM"
$issues.Range("E14").Value = "在同一个context中加载model，这样battery和programs.subprograms.testrecords.assignedbattery就是同一个instance."
$issues.Range("F14").Value = "Fixed"

# Row 14 now holds three lines of wrapped text, so it needs to grow taller.
$issues.Rows.Item(14).RowHeight = 108

# Re-activate the "Issues" sheet (it had lost focus to "Features") and move
# the selection there to F15.
$issues.Activate()
$issues.Range("F15").Select()
